$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.780.68'
$ws.Range("E2").Value = '  +2.62%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.765.48'
$ws.Range("E3").Value = '  -0.67%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.76%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.36'
$ws.Range("E5").Value = '  -0.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9970'
$ws.Range("E6").Value = '  -0.79%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3816'
$ws.Range("E7").Value = '  +0.57%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3423'
$ws.Range("E8").Value = '  +0.61%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.79'
$ws.Range("E9").Value = '  -2.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.136'
$ws.Range("E10").Value = '  -3.95%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07401'
$ws.Range("E11").Value = '  -0.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9973'
$ws.Range("E12").Value = '  -0.87%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.26'
$ws.Range("E13").Value = '  +3.56%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.344'
$ws.Range("E14").Value = '  -0.71%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.764.33'
$ws.Range("E15").Value = '  -0.77%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.069'
$ws.Range("E16").Value = '  +0.42%  '

$ws.Range("E17").Value = '  -0.69%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06663'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.15'
$ws.Range("E19").Value = '  -1.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9986'
$ws.Range("E20").Value = '  -0.51%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.32'
$ws.Range("E21").Value = '  +0.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.400'
$ws.Range("E22").Value = '  -2.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.783.11'
$ws.Range("E23").Value = '  +2.64%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.03'
$ws.Range("E24").Value = '  -1.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.374'
$ws.Range("E25").Value = '  -0.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.436'
$ws.Range("E26").Value = '  -1.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.68'
$ws.Range("E27").Value = '  -1.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.416'
$ws.Range("E28").Value = '  -3.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '153.18'
$ws.Range("E29").Value = '  -0.69%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.34'
$ws.Range("E30").Value = '  +0.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.962.42'
$ws.Range("E31").Value = '  -0.93%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.108'
$ws.Range("E32").Value = '  +2.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.948'
$ws.Range("E33").Value = '  -1.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08765'
$ws.Range("E34").Value = '  +1.46%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.74'
$ws.Range("E35").Value = '  -1.80%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02419'
$ws.Range("E36").Value = '  +4.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6792'
$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.311'
$ws.Range("E38").Value = '  -0.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06301'
$ws.Range("E39").Value = '  +0.59%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2180'
$ws.Range("E40").Value = '  +0.65%  '

$ws.Range("B41").Value = 'WEMIXTOKEN'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.503'
$ws.Range("E41").Value = '  -7.72%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.232'
$ws.Range("E42").Value = '  +0.50%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.228'
$ws.Range("E43").Value = '  -3.80%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.13'
$ws.Range("E44").Value = '  -0.38%  '

$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9984'
$ws.Range("E45").Value = '  -0.55%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6251'
$ws.Range("E46").Value = '  -1.78%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.835'
$ws.Range("E47").Value = '  -0.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.51'
$ws.Range("E48").Value = '  +0.75%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.077'
$ws.Range("E49").Value = '  -1.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07392'
$ws.Range("E50").Value = '  +4.16%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.95'
$ws.Range("E51").Value = '  -0.47%  '
